$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = "['HP:0012520', 'HP:0000709', 'HP:0002266', 'HP:0010651', 'HP:0020219', 'HP:0012379', 'HP:0007301', 'HP:0000098', 'HP:0002926', 'HP:0011195', 'HP:0001707', 'HP:0004323', 'HP:0011122', 'HP:0007598', 'HP:0010787', 'HP:0002286', 'HP:0100818', 'HP:0002121', 'HP:0001010', 'HP:0010719', 'HP:0000505', 'HP:0001533', 'HP:0000470', 'HP:0040195', 'HP:0001508', 'HP:0002373', 'HP:0000078', 'HP:0011360']"

$ws.Range("B12").Value = "['HP:0031938', 'HP:0011282', 'HP:0000582', 'HP:0000347', 'HP:0200043', 'HP:0001551', 'HP:0040195', 'HP:0025015', 'HP:0001762', 'HP:0003265', 'HP:0005164', 'HP:0000268', 'HP:0001933', 'HP:0011329', 'HP:0001647', 'HP:0001671']"

$ws.Range("B7").Select()
